# This script applies the Betfair odds update for the 2025-12-10 workbook
# (games / teams in rows 2-13 of Sheet1). Only numeric odds cells change;
# columns F:AO hold the back/lay odds, columns A:E (League/Date/Time/Home/Away)
# are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 1.49  # F2
$ws.Cells.Item(2, 7).Value = 1.51  # G2
$ws.Cells.Item(2, 8).Value = 7.4  # H2
$ws.Cells.Item(2, 9).Value = 8  # I2
$ws.Cells.Item(2, 10).Value = 4.9  # J2
$ws.Cells.Item(2, 11).Value = 5.2  # K2
$ws.Cells.Item(2, 15).Value = 1.21  # O2
$ws.Cells.Item(2, 16).Value = 2.44  # P2
$ws.Cells.Item(2, 17).Value = 1.65  # Q2
$ws.Cells.Item(2, 18).Value = 1.58  # R2
$ws.Cells.Item(2, 20).Value = 1.83  # T2
$ws.Cells.Item(2, 25).Value = 30  # Y2
$ws.Cells.Item(2, 26).Value = 70  # Z2
$ws.Cells.Item(2, 27).Value = 260  # AA2
$ws.Cells.Item(2, 28).Value = 11  # AB2
$ws.Cells.Item(2, 30).Value = 29  # AD2
$ws.Cells.Item(2, 31).Value = 120  # AE2
$ws.Cells.Item(2, 32).Value = 9.6  # AF2
$ws.Cells.Item(2, 33).Value = 10  # AG2
$ws.Cells.Item(2, 34).Value = 23  # AH2
$ws.Cells.Item(2, 36).Value = 13.5  # AJ2
$ws.Cells.Item(2, 37).Value = 16  # AK2
$ws.Cells.Item(2, 38).Value = 34  # AL2
$ws.Cells.Item(2, 39).Value = 130  # AM2
$ws.Cells.Item(2, 40).Value = 6  # AN2
$ws.Cells.Item(2, 41).Value = 130  # AO2

# Row 3
$ws.Cells.Item(3, 6).Value = 2  # F3
$ws.Cells.Item(3, 7).Value = 2.02  # G3
$ws.Cells.Item(3, 8).Value = 3.75  # H3
$ws.Cells.Item(3, 9).Value = 3.9  # I3
$ws.Cells.Item(3, 11).Value = 4.2  # K3
$ws.Cells.Item(3, 14).Value = 5.1  # N3
$ws.Cells.Item(3, 16).Value = 2.36  # P3
$ws.Cells.Item(3, 18).Value = 1.55  # R3
$ws.Cells.Item(3, 21).Value = 2.42  # U3
$ws.Cells.Item(3, 33).Value = 10.5  # AG3
$ws.Cells.Item(3, 37).Value = 24  # AK3
$ws.Cells.Item(3, 41).Value = 55  # AO3

# Row 5
$ws.Cells.Item(5, 7).Value = 1.6  # G5
$ws.Cells.Item(5, 8).Value = 1.09  # H5
$ws.Cells.Item(5, 16).Value = 1.96  # P5
$ws.Cells.Item(5, 17).Value = 1.6  # Q5

# Row 6
$ws.Cells.Item(6, 6).Value = 2.44  # F6
$ws.Cells.Item(6, 7).Value = 2.48  # G6
$ws.Cells.Item(6, 11).Value = 3.45  # K6
$ws.Cells.Item(6, 14).Value = 3.6  # N6
$ws.Cells.Item(6, 16).Value = 1.83  # P6
$ws.Cells.Item(6, 17).Value = 2.14  # Q6
$ws.Cells.Item(6, 18).Value = 1.32  # R6
$ws.Cells.Item(6, 20).Value = 1.86  # T6
$ws.Cells.Item(6, 21).Value = 2.08  # U6
$ws.Cells.Item(6, 31).Value = 44  # AE6
$ws.Cells.Item(6, 40).Value = 24  # AN6
$ws.Cells.Item(6, 41).Value = 46  # AO6

# Row 7
$ws.Cells.Item(7, 6).Value = 2.58  # F7
$ws.Cells.Item(7, 7).Value = 2.62  # G7
$ws.Cells.Item(7, 8).Value = 2.68  # H7
$ws.Cells.Item(7, 9).Value = 2.7  # I7
$ws.Cells.Item(7, 16).Value = 2.78  # P7
$ws.Cells.Item(7, 17).Value = 1.53  # Q7
$ws.Cells.Item(7, 24).Value = 29  # X7
$ws.Cells.Item(7, 25).Value = 18.5  # Y7
$ws.Cells.Item(7, 26).Value = 23  # Z7
$ws.Cells.Item(7, 27).Value = 42  # AA7
$ws.Cells.Item(7, 30).Value = 13  # AD7
$ws.Cells.Item(7, 31).Value = 25  # AE7
$ws.Cells.Item(7, 32).Value = 22  # AF7
$ws.Cells.Item(7, 36).Value = 40  # AJ7
$ws.Cells.Item(7, 40).Value = 13  # AN7
$ws.Cells.Item(7, 41).Value = 13.5  # AO7

# Row 8
$ws.Cells.Item(8, 7).Value = 10.5  # G8
$ws.Cells.Item(8, 11).Value = 5.8  # K8
$ws.Cells.Item(8, 16).Value = 2.5  # P8
$ws.Cells.Item(8, 20).Value = 1.97  # T8
$ws.Cells.Item(8, 24).Value = 24  # X8
$ws.Cells.Item(8, 28).Value = 38  # AB8
$ws.Cells.Item(8, 32).Value = 100  # AF8

# Row 9
$ws.Cells.Item(9, 9).Value = 22  # I9
$ws.Cells.Item(9, 11).Value = 8.4  # K9
$ws.Cells.Item(9, 18).Value = 1.76  # R9
$ws.Cells.Item(9, 31).Value = 440  # AE9
$ws.Cells.Item(9, 34).Value = 42  # AH9
$ws.Cells.Item(9, 35).Value = 290  # AI9
$ws.Cells.Item(9, 37).Value = 15  # AK9
$ws.Cells.Item(9, 38).Value = 1000  # AL9
$ws.Cells.Item(9, 40).Value = 3.35  # AN9

# Row 10
$ws.Cells.Item(10, 8).Value = 11.5  # H10
$ws.Cells.Item(10, 9).Value = 12.5  # I10
$ws.Cells.Item(10, 19).Value = 1.92  # S10
$ws.Cells.Item(10, 29).Value = 17  # AC10
$ws.Cells.Item(10, 33).Value = 11.5  # AG10
$ws.Cells.Item(10, 37).Value = 13  # AK10
$ws.Cells.Item(10, 38).Value = 28  # AL10
$ws.Cells.Item(10, 40).Value = 3.35  # AN10
$ws.Cells.Item(10, 41).Value = 1000  # AO10

# Row 11
$ws.Cells.Item(11, 10).Value = 4.3  # J11
$ws.Cells.Item(11, 16).Value = 2.14  # P11
$ws.Cells.Item(11, 25).Value = 9.2  # Y11

# Row 12
$ws.Cells.Item(12, 6).Value = 3.05  # F12
$ws.Cells.Item(12, 7).Value = 3.15  # G12
$ws.Cells.Item(12, 8).Value = 2.42  # H12
$ws.Cells.Item(12, 9).Value = 2.46  # I12
$ws.Cells.Item(12, 11).Value = 3.8  # K12
$ws.Cells.Item(12, 24).Value = 18.5  # X12
$ws.Cells.Item(12, 40).Value = 23  # AN12
$ws.Cells.Item(12, 41).Value = 16  # AO12

# Row 13
$ws.Cells.Item(13, 6).Value = 2.24  # F13
$ws.Cells.Item(13, 7).Value = 2.72  # G13
$ws.Cells.Item(13, 8).Value = 3.15  # H13
$ws.Cells.Item(13, 9).Value = 4.1  # I13
$ws.Cells.Item(13, 10).Value = 3.2  # J13
$ws.Cells.Item(13, 11).Value = 3.8  # K13
$ws.Cells.Item(13, 16).Value = 1.81  # P13
$ws.Cells.Item(13, 17).Value = 1.97  # Q13
